$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("BD2").Value = 151

# Row 3 updates
$ws.Range("G3").Value = 1.75
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 5.5
$ws.Range("L3").Value = 5.5
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("AJ3").Value = 51
$ws.Range("AR3").Value = 67
$ws.Range("AV3").Value = 67
